# Move the line-break between the "An Empty Plot" / "Project Description 1"
# / "Project Description 2" runs so that each later run starts with the
# newline instead of the previous run ending with it. The overall visible
# text is unchanged ("Figure D1.\nAn Empty Plot\nProject Description 1\n
# Project Description 2") -- only which run owns each newline changes.
#
# We edit the existing runs in place via TextRange.Characters(start, length)
# using the *original* character offsets, and we do so from the end of the
# text backwards so that earlier offsets are not invalidated by edits that
# change a run's length.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)        # "Title 16" placeholder holding the caption
$tr = $sh.TextFrame.TextRange

# Original run layout (1-based char offsets into the full text range):
#   Run1: start=1,  length=11  "Figure D1.\n"
#   Run2: start=12, length=14  "An Empty Plot\n"
#   Run3: start=26, length=22  "Project Description 1\n"
#   Run4: start=48, length=21  "Project Description 2"

# Run4: "Project Description 2" -> "\nProject Description 2"
$r4 = $tr.Characters(48, 21)
$r4.Text = "`nProject Description 2"

# Run3: "Project Description 1\n" -> "\nProject Description 1"
$r3 = $tr.Characters(26, 22)
$r3.Text = "`nProject Description 1"

# Run2: "An Empty Plot\n" -> "An Empty Plot"
$r2 = $tr.Characters(12, 14)
$r2.Text = "An Empty Plot"
